# Append three new daily rows (Aug 9-11, 2022) to the COVID tracking sheet,
# mirroring the format of the existing last row (row 34), and move the
# selection to the newly added block - matching the "data updated on Aug.12"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows' dates, as Excel serial date numbers (days since 1899-12-30).
# Row 34 (last existing row) holds 44781 = Aug 8 2022, so the new rows are
# Aug 9, 10 and 11 2022 -> 44782, 44783, 44784.
$newDates = @(44782, 44783, 44784)

$lastRow = 34

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $srcRow = $lastRow + $i
    $dstRow = $lastRow + $i + 1

    # Insert a new row, shifting down and inheriting the format of the row
    # above it (the row we are about to append after).
    $ws.Rows.Item($dstRow).Insert(-4121, 0) | Out-Null

    # The inserted row already carries over matching number formats/borders
    # for most columns; re-copy A, B, C and G individually from the row
    # above to make sure every column's style lines up with the template
    # row (single-cell copies apply cleanly here).
    $ws.Range("A" + $srcRow).Copy($ws.Range("A" + $dstRow))
    $ws.Range("B" + $srcRow).Copy($ws.Range("B" + $dstRow))
    $ws.Range("C" + $srcRow).Copy($ws.Range("C" + $dstRow))
    $ws.Range("G" + $srcRow).Copy($ws.Range("G" + $dstRow))

    # Fill in the row's values: the date in column A, and 0 (no new cases
    # reported) for every other tracked column, matching the template row.
    $ws.Range("A" + $dstRow).Value2 = $newDates[$i]
    $ws.Range("B" + $dstRow).Value2 = 0
    $ws.Range("C" + $dstRow).Value2 = 0
    $ws.Range("D" + $dstRow).Value2 = 0
    $ws.Range("E" + $dstRow).Value2 = 0
    $ws.Range("F" + $dstRow).Value2 = 0
    $ws.Range("G" + $dstRow).Value2 = 0
    $ws.Range("H" + $dstRow).Value2 = 0
    $ws.Range("I" + $dstRow).Value2 = 0
}

# Match the workbook's recorded selection after the edit: the newly added block.
$ws.Range("A35:I37").Select() | Out-Null
